# Update res_bus/vm_pu.xlsx: bus-2 slack voltage set-point changed from 1.05 pu to 1.02 pu
# (380 kV case) — every voltage-magnitude result in rows 2-25 shifts accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.021541367123396
$ws.Cells.Item(2, 4).Value2 = 1.026422218762465
$ws.Cells.Item(2, 5).Value2 = 1.025161480254848
$ws.Cells.Item(2, 6).Value2 = 1.019974042732294
$ws.Cells.Item(2, 9).Value2 = 1.028091854196729
$ws.Cells.Item(2, 10).Value2 = 1.026732252662637
$ws.Cells.Item(2, 11).Value2 = 1.029244911995427
$ws.Cells.Item(2, 12).Value2 = 1.027987861329245
$ws.Cells.Item(2, 13).Value2 = 1.022815702210854
$ws.Cells.Item(2, 14).Value2 = 1.028190330229254

# row 3
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.0227416582334
$ws.Cells.Item(3, 4).Value2 = 1.027302958779511
$ws.Cells.Item(3, 5).Value2 = 1.026309295284918
$ws.Cells.Item(3, 6).Value2 = 1.021817059437145
$ws.Cells.Item(3, 9).Value2 = 1.028285991724038
$ws.Cells.Item(3, 10).Value2 = 1.027568902122241
$ws.Cells.Item(3, 11).Value2 = 1.029932964539345
$ws.Cells.Item(3, 12).Value2 = 1.02894199585776
$ws.Cells.Item(3, 13).Value2 = 1.024462013769011
$ws.Cells.Item(3, 14).Value2 = 1.029028167827055

# row 4
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.023517449535878
$ws.Cells.Item(4, 4).Value2 = 1.027871857132009
$ws.Cells.Item(4, 5).Value2 = 1.027051514385379
$ws.Cells.Item(4, 6).Value2 = 1.023008430115822
$ws.Cells.Item(4, 9).Value2 = 1.028409832110904
$ws.Cells.Item(4, 10).Value2 = 1.02810892762239
$ws.Cells.Item(4, 11).Value2 = 1.03037655236141
$ws.Cells.Item(4, 12).Value2 = 1.029558321532134
$ws.Cells.Item(4, 13).Value2 = 1.025525699315718
$ws.Cells.Item(4, 14).Value2 = 1.029568960225357

# row 5
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.02384338607786
$ws.Cells.Item(5, 4).Value2 = 1.028110784678137
$ws.Cells.Item(5, 5).Value2 = 1.027363428223793
$ws.Cells.Item(5, 6).Value2 = 1.023509010408098
$ws.Cells.Item(5, 9).Value2 = 1.028461469200465
$ws.Cells.Item(5, 10).Value2 = 1.028335635680955
$ws.Cells.Item(5, 11).Value2 = 1.030562648853155
$ws.Cells.Item(5, 12).Value2 = 1.029817173258698
$ws.Cells.Item(5, 13).Value2 = 1.025972502796357
$ws.Cells.Item(5, 14).Value2 = 1.029795990235368

# row 6
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.023898100283504
$ws.Cells.Item(6, 4).Value2 = 1.028150887803256
$ws.Cells.Item(6, 5).Value2 = 1.027415793274114
$ws.Cells.Item(6, 6).Value2 = 1.023593044430947
$ws.Cells.Item(6, 9).Value2 = 1.028470114366826
$ws.Cells.Item(6, 10).Value2 = 1.028373682355189
$ws.Cells.Item(6, 11).Value2 = 1.030593872587168
$ws.Cells.Item(6, 12).Value2 = 1.029860620934364
$ws.Cells.Item(6, 13).Value2 = 1.026047501671828
$ws.Cells.Item(6, 14).Value2 = 1.029834090940244

# row 7
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.023521805520994
$ws.Cells.Item(7, 4).Value2 = 1.027875050623741
$ws.Cells.Item(7, 5).Value2 = 1.027055682641579
$ws.Cells.Item(7, 6).Value2 = 1.023015119948415
$ws.Cells.Item(7, 9).Value2 = 1.028410523759031
$ws.Cells.Item(7, 10).Value2 = 1.028111958155284
$ws.Cells.Item(7, 11).Value2 = 1.030379040513347
$ws.Cells.Item(7, 12).Value2 = 1.029561781306526
$ws.Cells.Item(7, 13).Value2 = 1.02553167096657
$ws.Cells.Item(7, 14).Value2 = 1.029571995061955

# row 8
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.021947195134762
$ws.Cells.Item(8, 4).Value2 = 1.026720077055722
$ws.Cells.Item(8, 5).Value2 = 1.025549493924741
$ws.Cells.Item(8, 6).Value2 = 1.020597151087282
$ws.Cells.Item(8, 9).Value2 = 1.028157832734356
$ws.Cells.Item(8, 10).Value2 = 1.02701528194199
$ws.Cells.Item(8, 11).Value2 = 1.029477780886488
$ws.Cells.Item(8, 12).Value2 = 1.028310537470583
$ws.Cells.Item(8, 13).Value2 = 1.023372415940752
$ws.Cells.Item(8, 14).Value2 = 1.028473761442647

# row 9
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.019165635362018
$ws.Cells.Item(9, 4).Value2 = 1.02467712513055
$ws.Cells.Item(9, 5).Value2 = 1.022891469864513
$ws.Cells.Item(9, 6).Value2 = 1.016326789294316
$ws.Cells.Item(9, 9).Value2 = 1.027698900929648
$ws.Cells.Item(9, 10).Value2 = 1.025072390231494
$ws.Cells.Item(9, 11).Value2 = 1.027877087460865
$ws.Cells.Item(9, 12).Value2 = 1.02609741080653
$ws.Cells.Item(9, 13).Value2 = 1.019554908111893
$ws.Cells.Item(9, 14).Value2 = 1.026528110603069

# row 10
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.017306377644079
$ws.Cells.Item(10, 4).Value2 = 1.023309827395275
$ws.Cells.Item(10, 5).Value2 = 1.021116609445234
$ws.Cells.Item(10, 6).Value2 = 1.013472702899782
$ws.Cells.Item(10, 9).Value2 = 1.027383722235323
$ws.Cells.Item(10, 10).Value2 = 1.023769953965848
$ws.Cells.Item(10, 11).Value2 = 1.02680138773708
$ws.Cells.Item(10, 12).Value2 = 1.024616241269252
$ws.Cells.Item(10, 13).Value2 = 1.017000768079592
$ws.Cells.Item(10, 14).Value2 = 1.025223824728534

# row 11
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.016500083846119
$ws.Cells.Item(11, 4).Value2 = 1.022716479371972
$ws.Cells.Item(11, 5).Value2 = 1.020347354358444
$ws.Cells.Item(11, 6).Value2 = 1.012234986632489
$ws.Cells.Item(11, 9).Value2 = 1.027245048871187
$ws.Cells.Item(11, 10).Value2 = 1.023204242631317
$ws.Cells.Item(11, 11).Value2 = 1.026333535939186
$ws.Cells.Item(11, 12).Value2 = 1.023973471620808
$ws.Cells.Item(11, 13).Value2 = 1.015892486272606
$ws.Cells.Item(11, 14).Value2 = 1.02465731001902

# row 12
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.016200401763875
$ws.Cells.Item(12, 4).Value2 = 1.0224958857279
$ws.Cells.Item(12, 5).Value2 = 1.020061505186056
$ws.Cells.Item(12, 6).Value2 = 1.011774947486508
$ws.Cells.Item(12, 9).Value2 = 1.027193208149032
$ws.Cells.Item(12, 10).Value2 = 1.022993846103624
$ws.Cells.Item(12, 11).Value2 = 1.02615944184356
$ws.Cells.Item(12, 12).Value2 = 1.023734502226764
$ws.Cells.Item(12, 13).Value2 = 1.015480459670389
$ws.Cells.Item(12, 14).Value2 = 1.024446614704125

# row 13
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.016264693204942
$ws.Cells.Item(13, 4).Value2 = 1.022543212804886
$ws.Cells.Item(13, 5).Value2 = 1.0201228260097
$ws.Cells.Item(13, 6).Value2 = 1.011873641076299
$ws.Cells.Item(13, 9).Value2 = 1.02720434316469
$ws.Cells.Item(13, 10).Value2 = 1.023038988974317
$ws.Cells.Item(13, 11).Value2 = 1.026196799828921
$ws.Cells.Item(13, 12).Value2 = 1.023785771814102
$ws.Cells.Item(13, 13).Value2 = 1.015568857330613
$ws.Cells.Item(13, 14).Value2 = 1.024491821682872

# row 14
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.016475315909108
$ws.Cells.Item(14, 4).Value2 = 1.02269824908211
$ws.Cells.Item(14, 5).Value2 = 1.020323728317493
$ws.Cells.Item(14, 6).Value2 = 1.012196965777625
$ws.Cells.Item(14, 9).Value2 = 1.027240770463772
$ws.Cells.Item(14, 10).Value2 = 1.023186856638773
$ws.Cells.Item(14, 11).Value2 = 1.026319151669391
$ws.Cells.Item(14, 12).Value2 = 1.023953722767014
$ws.Cells.Item(14, 13).Value2 = 1.015858435470055
$ws.Cells.Item(14, 14).Value2 = 1.024639899336373

# row 15
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.016605062360325
$ws.Cells.Item(15, 4).Value2 = 1.022793745764756
$ws.Cells.Item(15, 5).Value2 = 1.020447495686334
$ws.Cells.Item(15, 6).Value2 = 1.012396137123163
$ws.Cells.Item(15, 9).Value2 = 1.027263170602699
$ws.Cells.Item(15, 10).Value2 = 1.023277927412693
$ws.Cells.Item(15, 11).Value2 = 1.026394495108303
$ws.Cells.Item(15, 12).Value2 = 1.024057174139837
$ws.Cells.Item(15, 13).Value2 = 1.016036805855596
$ws.Cells.Item(15, 14).Value2 = 1.024731099441237

# row 16
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.017359862521946
$ws.Cells.Item(16, 4).Value2 = 1.023349178358322
$ws.Cells.Item(16, 5).Value2 = 1.021167646668684
$ws.Cells.Item(16, 6).Value2 = 1.013554805286523
$ws.Cells.Item(16, 9).Value2 = 1.027392879124819
$ws.Cells.Item(16, 10).Value2 = 1.023807461236221
$ws.Cells.Item(16, 11).Value2 = 1.026832393754919
$ws.Cells.Item(16, 12).Value2 = 1.024658869672815
$ws.Cells.Item(16, 13).Value2 = 1.017074271100631
$ws.Cells.Item(16, 14).Value2 = 1.025261385263533

# row 17
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.017832997836075
$ws.Cells.Item(17, 4).Value2 = 1.02369723692485
$ws.Cells.Item(17, 5).Value2 = 1.021619180356294
$ws.Cells.Item(17, 6).Value2 = 1.014281094520243
$ws.Cells.Item(17, 9).Value2 = 1.027473652512086
$ws.Cells.Item(17, 10).Value2 = 1.024139153366864
$ws.Cells.Item(17, 11).Value2 = 1.027106520856449
$ws.Cells.Item(17, 12).Value2 = 1.025035916308233
$ws.Cells.Item(17, 13).Value2 = 1.017724415101677
$ws.Cells.Item(17, 14).Value2 = 1.025593548435047

# row 18
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.018108852000432
$ws.Cells.Item(18, 4).Value2 = 1.023900128433563
$ws.Cells.Item(18, 5).Value2 = 1.021882482275929
$ws.Cells.Item(18, 6).Value2 = 1.014704546496055
$ws.Cells.Item(18, 9).Value2 = 1.027520554161131
$ws.Cells.Item(18, 10).Value2 = 1.024332455376692
$ws.Cells.Item(18, 11).Value2 = 1.027266215326748
$ws.Cells.Item(18, 12).Value2 = 1.025255705109996
$ws.Cells.Item(18, 13).Value2 = 1.018103409876865
$ws.Cells.Item(18, 14).Value2 = 1.025787124955902

# row 19
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.018202891332944
$ws.Cells.Item(19, 4).Value2 = 1.0239692880425
$ws.Cells.Item(19, 5).Value2 = 1.021972249668697
$ws.Cells.Item(19, 6).Value2 = 1.014848902527553
$ws.Cells.Item(19, 9).Value2 = 1.027536510462235
$ws.Cells.Item(19, 10).Value2 = 1.024398337910228
$ws.Cells.Item(19, 11).Value2 = 1.027320633299444
$ws.Cells.Item(19, 12).Value2 = 1.025330624436169
$ws.Cells.Item(19, 13).Value2 = 1.018232599815175
$ws.Cells.Item(19, 14).Value2 = 1.025853101050192

# row 20
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.017782247101731
$ws.Cells.Item(20, 4).Value2 = 1.023659906500768
$ws.Cells.Item(20, 5).Value2 = 1.021570742340366
$ws.Cells.Item(20, 6).Value2 = 1.01420318930705
$ws.Cells.Item(20, 9).Value2 = 1.027465008237561
$ws.Cells.Item(20, 10).Value2 = 1.024103583385342
$ws.Cells.Item(20, 11).Value2 = 1.027077130239813
$ws.Cells.Item(20, 12).Value2 = 1.02499547688068
$ws.Cells.Item(20, 13).Value2 = 1.017654683984797
$ws.Cells.Item(20, 14).Value2 = 1.025557927940072

# row 21
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.016413298020266
$ws.Cells.Item(21, 4).Value2 = 1.02265260024306
$ws.Cells.Item(21, 5).Value2 = 1.020264570763629
$ws.Cells.Item(21, 6).Value2 = 1.012101763002226
$ws.Cells.Item(21, 9).Value2 = 1.0272300526838
$ws.Cells.Item(21, 10).Value2 = 1.023143320672861
$ws.Cells.Item(21, 11).Value2 = 1.026283130766257
$ws.Cells.Item(21, 12).Value2 = 1.023904271394071
$ws.Cells.Item(21, 13).Value2 = 1.015773172002556
$ws.Cells.Item(21, 14).Value2 = 1.024596301544396

# row 22
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.015551491367961
$ws.Cells.Item(22, 4).Value2 = 1.022018121277493
$ws.Cells.Item(22, 5).Value2 = 1.019442669650263
$ws.Cells.Item(22, 6).Value2 = 1.010778792024044
$ws.Cells.Item(22, 9).Value2 = 1.027080409778545
$ws.Cells.Item(22, 10).Value2 = 1.0225380237471
$ws.Cells.Item(22, 11).Value2 = 1.025782098933387
$ws.Cells.Item(22, 12).Value2 = 1.023216934856809
$ws.Cells.Item(22, 13).Value2 = 1.014588093576365
$ws.Cells.Item(22, 14).Value2 = 1.023990145027573

# row 23
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.016008456738315
$ws.Cells.Item(23, 4).Value2 = 1.022354580040661
$ws.Cells.Item(23, 5).Value2 = 1.019878438839293
$ws.Cells.Item(23, 6).Value2 = 1.01148029181217
$ws.Cells.Item(23, 9).Value2 = 1.027159920328494
$ws.Cells.Item(23, 10).Value2 = 1.022859050425892
$ws.Cells.Item(23, 11).Value2 = 1.026047878101897
$ws.Cells.Item(23, 12).Value2 = 1.023581424934446
$ws.Cells.Item(23, 13).Value2 = 1.015216529055328
$ws.Cells.Item(23, 14).Value2 = 1.02431162760107

# row 24
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.017805179537282
$ws.Cells.Item(24, 4).Value2 = 1.023676774899179
$ws.Cells.Item(24, 5).Value2 = 1.021592629610964
$ws.Cells.Item(24, 6).Value2 = 1.014238391875079
$ws.Cells.Item(24, 9).Value2 = 1.027468914868453
$ws.Cells.Item(24, 10).Value2 = 1.024119656448996
$ws.Cells.Item(24, 11).Value2 = 1.027090411209816
$ws.Cells.Item(24, 12).Value2 = 1.025013750137663
$ws.Cells.Item(24, 13).Value2 = 1.017686193163495
$ws.Cells.Item(24, 14).Value2 = 1.02557402382932

# row 25
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.019885577824228
$ws.Cells.Item(25, 4).Value2 = 1.025206207304648
$ws.Cells.Item(25, 5).Value2 = 1.023579119069293
$ws.Cells.Item(25, 6).Value2 = 1.017431991640684
$ws.Cells.Item(25, 9).Value2 = 1.027819168555196
$ws.Cells.Item(25, 10).Value2 = 1.025575925876042
$ws.Cells.Item(25, 11).Value2 = 1.028292406342592
$ws.Cells.Item(25, 12).Value2 = 1.026670556996481
$ws.Cells.Item(25, 13).Value2 = 1.020543387769453
$ws.Cells.Item(25, 14).Value2 = 1.027032361325988
